$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB24 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$newB25 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"

$ws.Range("B24").Value = $newB24
$ws.Range("C24").Value = $newB24
$ws.Range("B25").Value = $newB25
$ws.Range("C25").Value = $newB25
